$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "1" to "მარნეული"
$ws.Name = "მარნეული"

# Remove the obsolete "(მოსახლეობის აღწერის შედეგებით)" row (old row 2) -
# rows below shift up by one.
$ws.Rows(2).Delete()

# Drop the 1989 / 2002 columns, leaving only the 2014 figures (old columns
# B and C collapse away, old column D becomes the new column B).
$ws.Range("B:C").Delete()

# Match the saved selection/active-cell state.
$null = $ws.Range("A2").Select()
